$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "Dittmann"

# Row 17: Kalenderwoche 22 -> 23, add "Fertigstellung der Dokumentation" task, 2 hours
$ws.Range("B17").Value = 23
$ws.Range("C17").Value = "Fertigstellung der Dokumentation"
$ws.Range("C17").HorizontalAlignment = -4108
$ws.Range("C17").VerticalAlignment = -4108
$ws.Range("C17").WrapText = $true
$ws.Range("D17").Value = 2

# Row 18: Kalenderwoche 23 -> 26, add "Fertigstellung der Dokumentation" task, 2 hours
$ws.Range("B18").Value = 26
$ws.Range("C18").Value = "Fertigstellung der Dokumentation"
$ws.Range("C18").HorizontalAlignment = -4108
$ws.Range("C18").VerticalAlignment = -4108
$ws.Range("C18").WrapText = $true
$ws.Range("D18").Value = 2

# Switch active sheet to "Dittmann" and update selection
$ws.Activate()
$ws.Range("G16").Select()
